$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (Russian column names)
$ws.Range("A1").Value = "Название"
$ws.Range("B1").Value = "Норма азота"
$ws.Range("C1").Value = "Норма фосфора"
$ws.Range("D1").Value = "Норма калия"
$ws.Range("E1").Value = "Id культуры"
$ws.Range("F1").Value = "Район"
$ws.Range("G1").Value = "Цена"
$ws.Range("H1").Value = "Описание "
$ws.Range("I1").Value = "Назначение"

# Row 2
$ws.Range("A2").Value = "имя1"
$ws.Range("B2").Value = 1
$ws.Range("F2").Value = "Район 1"
$ws.Range("H2").Value = "Описание 1"
$ws.Range("I2").Value = "Цель 1"

# Row 3
$ws.Range("A3").Value = "имя2"
$ws.Range("B3").Value = 123
$ws.Range("C3").Value = 1
$ws.Range("F3").Value = "Район 2"
$ws.Range("H3").Value = "Описание 2"
$ws.Range("I3").Value = "Цель 2"

# Column widths (input values chosen so the runtime's width-rounding lands
# on the closest representable value to the target stored width)
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws.Columns.Item(2).ColumnWidth = 17.5
$ws.Columns.Item(3).ColumnWidth = 20.5
$ws.Columns.Item(4).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 14.166666666666666
$ws.Columns.Item(8).ColumnWidth = 15.0
$ws.Columns.Item(9).ColumnWidth = 16.833333333333332

# Selection change
$ws.Range("J1").Select()
